# Update workbook to roll the SIPOT report from Q3-2021 to Q4-2021.
# - Moves the sheet view so A3 is the top-left visible cell and A8 is selected
#   (previously the view was scrolled to E3 with G3:I3 selected).
# - Updates the period-covered dates (B8/C8) and the report/validation dates
#   (AB8/AC8) to the following quarter.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")
$ws.Activate()

# --- Update the date values for row 8 (keep existing number formatting/style) ---
$ws.Range("B8").Value = 44470   # 2021-10-01 (was 44378 / 2021-07-01)
$ws.Range("C8").Value = 44561   # 2021-12-31 (was 44469 / 2021-09-30)
$ws.Range("AB8").Value = 44571  # 2022-01-10 (was 44480 / 2021-10-11)
$ws.Range("AC8").Value = 44571  # 2022-01-10 (was 44480 / 2021-10-11)

# --- Update the sheet view / selection state ---
# Scroll so column A / row 3 become the top-left visible cell (was E3).
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 3

# Select A8 as the active cell (was G3:I3 with G3 active).
$ws.Range("A8").Select()
